$d = $word.ActiveDocument
for ($f = 1; $f -le 3; $f++) {
    try {
        $ftr = $d.Sections.Item(1).Footers.Item($f)
        if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
            $shp = $ftr.Range.InlineShapes.Item(1)
            $shp.Name = "TESTNAME"
            Write-Output ("Footer " + $f + ": set OK, Name=" + $shp.Name)
        } else {
            Write-Output ("Footer " + $f + ": no shapes or doesn't exist")
        }
    } catch {
        Write-Output ("Footer " + $f + ": EXCEPTION " + $_)
    }
}
for ($h = 1; $h -le 3; $h++) {
    try {
        $hdr = $d.Sections.Item(1).Headers.Item($h)
        if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
            $shp = $hdr.Range.InlineShapes.Item(1)
            $shp.Name = "TESTNAME"
            Write-Output ("Header " + $h + ": set OK, Name=" + $shp.Name)
        } else {
            Write-Output ("Header " + $h + ": no shapes or doesn't exist")
        }
    } catch {
        Write-Output ("Header " + $h + ": EXCEPTION " + $_)
    }
}
